$d = $word.ActiveDocument

# Insert the contact-info line as a new centered paragraph immediately after
# the "Dheeraj Chand" name line, using a Find/Replace with a paragraph mark
# (^p) so the new paragraph naturally inherits the centered alignment of the
# name paragraph without picking up the name run's bold/size character
# formatting.
$d.Content.Find.Execute(
    "Dheeraj Chand", $false, $false, $false, $false, $false,
    $true, 1, $false,
    "Dheeraj Chand^p202.550.7110 | dheeraj.chand@gmail.com | https://www.dheerajchand.com | https://www.linkedin.com/in/dheerajchand/ | Austin, TX",
    2
)
